$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mike took on the grid-search / hyper-parameter tuning task (row 3 of the
# "Final Project Work" resource table) -- assign the resource in F3.
$ws.Range("F3").Value = "Mike"

# Move the active selection to F4, matching the author's cursor position
# after making the edit.
$ws.Range("F4").Select()
